$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3092107474803925
$ws.Range("B1").Value = 0.333082526922226
$ws.Range("C1").Value = 5.513103485107422
$ws.Range("D1").Value = 2.831963062286377
$ws.Range("E1").Value = 1.282640218734741
